$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update well labels from F1..F12 to G1..G12 in column A (rows 2-13)
for ($i = 1; $i -le 12; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = "G$i"
}

# Update the selected range on the sheet view
$ws.Range("A2:A13").Select()
